$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the old table (content + formatting) ---
$ws.Cells.Clear()

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 11.140625
$ws.Columns.Item(2).ColumnWidth = 9.42578125
$ws.Columns.Item(3).ColumnWidth = 25.42578125
$ws.Columns.Item(4).ColumnWidth = 25.85546875
$ws.Columns.Item(5).ColumnWidth = 27.28515625
$ws.Columns.Item(6).ColumnWidth = 33.140625
$ws.Columns.Item(7).ColumnWidth = 38

# --- Header row ---
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Hora"
$ws.Range("C1").Value = "WC47 NACP"
$ws.Range("D1").Value = "WC48 POWER 5F"
$ws.Range("E1").Value = "WC49 POWER 5H"
$ws.Range("F1").Value = "WV50 FILTER"
$ws.Range("G1").Value = "SPL"

# Re-apply the header style (yellow fill, centered, bordered) that already
# exists in the workbook (style index 1), then thin the border out.
$ws.Range("A1:G1").Style = $ws.Range("A1").Style
$hdr = $ws.Range("A1:F1")
$hdr.Interior.Pattern = 1
$hdr.Interior.PatternColorIndex = -4105
$hdr.Interior.Color = 65535
$hdr.HorizontalAlignment = -4108
$ws.Range("A1:G1").Interior.Color = 65535
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").Borders.LineStyle = 1
$ws.Range("A1:G1").Borders.Weight = 2

# --- Data rows ---
$ws.Range("A2").Value = "'2024-05-10"
$ws.Range("B2").Value = "11:54:47"
$ws.Range("C2").Value = "Palet atascado en la curva"

$ws.Range("A3").Value = "'2024-05-10"
$ws.Range("B3").Value = "11:54:49"
$ws.Range("F3").Value = "QR desplazado"

$ws.Range("A4").Value = "'2024-05-10"
$ws.Range("B4").Value = "11:54:51"
$ws.Range("G4").Value = "Error en sensor de salida"

$ws.Range("A5").Value = "'2024-05-10"
$ws.Range("B5").Value = "11:54:54"
$ws.Range("D5").Value = "Fallo etiqueta"

# --- Row 1 should go back to the default (non-custom) row height ---
$ws.Rows.Item(1).EntireRow.AutoFit()

# --- Selection, matching the saved cursor position ---
$ws.Range("C7:C8").Select()

Write-Output "edit applied"
